{"js": "// 1) Merge the split \"thousands\" sentence back into a single run and fix\n//    the typo: \"thousands rubles\" -> \"thousand rubles\" (removing the\n//    proofing-error markers around \"thousands\" in the process).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet pricingParagraph = null;\nlet conclusionParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const t = p.text;\n  if (pricingParagraph === null && t.indexOf(\"As mentioned above\") === 0) {\n    pricingParagraph = p;\n  }\n  if (conclusionParagraph === null && t.trim() === \"Conclusion\") {\n    conclusionParagraph = p;\n  }\n}\n\nif (!pricingParagraph) {\n  throw new Error(\"Could not find the 'As mentioned above...' paragraph.\");\n}\nif (!conclusionParagraph) {\n  throw new Error(\"Could not find the 'Conclusion' heading paragraph.\");\n}\n\npricingParagraph.load(\"text\");\nawait context.sync();\n\nconst fixedText = pricingParagraph.text.replace(\n  \"six and a half thousands rubles\",\n  \"six and a half thousand rubles\"\n);\n\n// Re-insert the whole paragraph text as plain text. This collapses the\n// previous three runs (split around the <w:proofErr> gramStart/gramEnd\n// markers wrapping \"thousands\") into a single run and drops those\n// proofing-error markers, matching the target document structure.\npricingParagraph.getRange().insertText(fixedText, \"Replace\");\n\n// 2) Add a hanging indent (0.18in / 216 twips) to the \"Conclusion\"\n//    Heading 1 paragraph, matching the indent already implied by the\n//    Heading1 style's tab stop.\nconclusionParagraph.leftIndent = 10.8;\nconclusionParagraph.firstLineIndent = -10.8;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the pricing paragraph (\"As mentioned above...\") and the\n# \"Conclusion\" Heading 1 paragraph by their text, rather than a hard-coded\n# index, so the script is resilient to unrelated paragraph shifts.\n$pricingParagraph = $null\n$conclusionParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($pricingParagraph -eq $null -and $t.StartsWith(\"As mentioned above\")) {\n        $pricingParagraph = $p\n    }\n    if ($conclusionParagraph -eq $null -and $t.Trim() -eq \"Conclusion\") {\n        $conclusionParagraph = $p\n    }\n}\n\nif ($pricingParagraph -eq $null) {\n    throw \"Could not find the 'As mentioned above...' paragraph.\"\n}\nif ($conclusionParagraph -eq $null) {\n    throw \"Could not find the 'Conclusion' heading paragraph.\"\n}\n\n# 1) Fix the typo \"thousands rubles\" -> \"thousand rubles\". Doing this via\n#    Find/Replace on the paragraph's range re-writes the paragraph as a\n#    single run, which also collapses the <w:proofErr> gramStart/gramEnd\n#    markers that previously wrapped the word \"thousands\".\n$find = $pricingParagraph.Range.Find\n$find.Text = \"six and a half thousands rubles\"\n$find.Replacement.Text = \"six and a half thousand rubles\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# 2) Add a hanging indent (0.18in / 216 twips) to the \"Conclusion\" heading,\n#    matching the indent implied by the Heading1 style's tab stop.\n$conclusionParagraph.Format.LeftIndent = 10.8\n$conclusionParagraph.Format.FirstLineIndent = -10.8\n"}
